# Apply the "started editing student responses" commit to the codebook workbook.
#
# Net effect (once the shared-string bookkeeping is accounted for):
#   1. On the "educator" sheet, the code in A7 changes from "learning_interests"
#      to "teaching_interests" (a new distinct code label is introduced).
#   2. On the "educator" sheet, the description in B27 (for the
#      "student_enjoyment" code) has a typo fixed: "how will did your
#      students" -> "how well did your students".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("educator")

$ws.Range("A7").Value = "teaching_interests"

$ws.Range("B27").Value = "Compared to how you taught similar material previously,  how well did your students seem to enjoy  the case studies?"
